$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Move the existing "total" row (17) and "footer" row (18) down to make
#    room for two new data rows (the table grows from 10 to 12 items).
#    Copy formats+values from row 18 -> row 20, and row 17 -> row 19.
# ---------------------------------------------------------------------------
$ws.Range("A18:Q18").Copy()
$ws.Range("A20:Q20").PasteSpecial(-4104)

$ws.Range("A17:Q17").Copy()
$ws.Range("A19:Q19").PasteSpecial(-4104)
$ws.Application.CutCopyMode = 0

$ws.Rows("20:20").RowHeight = 16.5
$ws.Rows("19:19").RowHeight = 25.5

# Merges for the relocated rows
$ws.Range("P19:Q19").Merge()
$ws.Range("A20:F20").Merge()
$ws.Range("G20:I20").Merge()
$ws.Range("K20:Q20").Merge()

# Update the grand total (sum of "sale price" column for all 12 items) and
# the generation timestamp in the footer.
$ws.Range("P19").Value = 1085.24
$ws.Range("A20").Value = "Friday, 15 August, 2025 6:58 PM"

# ---------------------------------------------------------------------------
# 2) Turn (old) rows 17 & 18 into two brand-new data rows, matching the
#    look (style/format/merges) of the existing item rows.
# ---------------------------------------------------------------------------
$ws.Range("A16:Q16").Copy()
$ws.Range("A17:Q17").PasteSpecial(-4104)
$ws.Range("A16:Q16").Copy()
$ws.Range("A18:Q18").PasteSpecial(-4104)
$ws.Application.CutCopyMode = 0

$ws.Rows("17:17").RowHeight = 25.5
$ws.Rows("18:18").RowHeight = 24.75

$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

$ws.Range("A18:B18").Merge()
$ws.Range("C18:G18").Merge()
$ws.Range("H18:K18").Merge()
$ws.Range("L18:M18").Merge()
$ws.Range("N18:O18").Merge()

function Set-TextCell($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
}

# Row 17 -> "سرنجات 3 سم" (was item #9, now item #11)
$ws.Range("A17").Value = 11
Set-TextCell "C17" "سرنجات 3 سم"
Set-TextCell "H17" "0:0"
Set-TextCell "L17" "0"
Set-TextCell "N17" "2.00"
Set-TextCell "P17" "12.0000"
Set-TextCell "Q17" "6:0"

# Row 18 -> "كالونا" (was item #10, now item #12)
$ws.Range("A18").Value = 12
Set-TextCell "C18" "كالونا "
Set-TextCell "H18" "0:0"
Set-TextCell "L18" "0"
Set-TextCell "N18" "15.00"
Set-TextCell "P18" "15.0000"
Set-TextCell "Q18" "1:0"

# ---------------------------------------------------------------------------
# 3) Re-order the original 10 item rows (7-16) so the two newly received
#    products - CIPROFAR and TRESIBA - are inserted alphabetically, shifting
#    the rest of the list down by one each.
# ---------------------------------------------------------------------------

# Row 7: CATAFLAM (unchanged)
Set-TextCell "C7" "CATAFLAM 75MG/3ML 6 AMP."
Set-TextCell "H7" "1:0"
Set-TextCell "L7" "1"
Set-TextCell "N7" "120.00"
Set-TextCell "P7" "19.2000"
Set-TextCell "Q7" "0:1"

# Row 8: CIPROFAR (new item)
Set-TextCell "C8" "CIPROFAR 500MG 10 F.C.TAB"
Set-TextCell "H8" "0:0"
Set-TextCell "L8" "1"
Set-TextCell "N8" "69.00"
Set-TextCell "P8" "69.0000"
Set-TextCell "Q8" "1:0"

# Row 9: FATROXIM (was row 8)
Set-TextCell "C9" "FATROXIM 550 MG 30TAB"
Set-TextCell "H9" "0:0"
Set-TextCell "L9" "0"
Set-TextCell "N9" "360.00"
Set-TextCell "P9" "118.8000"
Set-TextCell "Q9" "0:1"

# Row 10: KETOLAC (was row 9)
Set-TextCell "C10" "KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF."
Set-TextCell "H10" "3:0"
Set-TextCell "L10" "1"
Set-TextCell "N10" "60.00"
Set-TextCell "P10" "36.0000"
Set-TextCell "Q10" "0:3"

# Row 11: MOBITIL (was row 10)
Set-TextCell "C11" "MOBITIL 15MG/1.5ML 3 AMP."
Set-TextCell "H11" "2:2"
Set-TextCell "L11" "1"
Set-TextCell "N11" "39.00"
Set-TextCell "P11" "25.7400"
Set-TextCell "Q11" "0:2"

# Row 12: NORHINOSE (was row 11)
Set-TextCell "C12" "NORHINOSE 50MCG/DOSE NASAL SPRAY 120 DOSES"
Set-TextCell "H12" "3:0"
Set-TextCell "L12" "1"
Set-TextCell "N12" "90.00"
Set-TextCell "P12" "90.0000"
Set-TextCell "Q12" "1:0"

# Row 13: PAROFEN (was row 12)
Set-TextCell "C13" "PAROFEN 30 OBLONG TAB."
Set-TextCell "H13" "0:2"
Set-TextCell "L13" "1"
Set-TextCell "N13" "69.00"
Set-TextCell "P13" "69.0000"
Set-TextCell "Q13" "1:0"

# Row 14: TRESIBA (new item, was URGINAFECT's slot)
Set-TextCell "C14" "TRESIBA 100 I.U./ML FLEXTOUCH PRE-FILLED PEN"
Set-TextCell "H14" "0:0"
Set-TextCell "L14" "1"
Set-TextCell "N14" "550.00"
Set-TextCell "P14" "550.0000"
Set-TextCell "Q14" "1:0"

# Row 15: URGINAFECT (was row 13)
Set-TextCell "C15" "URGINAFECT 10MG 20 F.C. TAB"
Set-TextCell "H15" "0:1"
Set-TextCell "L15" "1"
Set-TextCell "N15" "71.00"
Set-TextCell "P15" "35.5000"
Set-TextCell "Q15" "0:1"

# Row 16: زيت فاتيكا (was row 14)
Set-TextCell "C16" "زيت فاتيكا كبير 180 مل"
Set-TextCell "H16" "12:0"
Set-TextCell "L16" "0"
Set-TextCell "N16" "45.00"
Set-TextCell "P16" "45.0000"
Set-TextCell "Q16" "1:0"

Write-Output "edit complete"
